$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrBE = New-Object "object[,]" 24,4
$arrBE[0,0] = 18.91363166473712
$arrBE[0,1] = 15.28701316329902
$arrBE[0,2] = 15.03032962933279
$arrBE[0,3] = 16.4524776792096
$arrBE[1,0] = 18.39786324980298
$arrBE[1,1] = 14.7895303551205
$arrBE[1,2] = 14.97548227606261
$arrBE[1,3] = 16.39854111648452
$arrBE[2,0] = 18.07860768532085
$arrBE[2,1] = 14.47978653323477
$arrBE[2,2] = 14.94543175829169
$arrBE[2,3] = 16.36941931823545
$arrBE[3,0] = 17.94807296541402
$arrBE[3,1] = 14.35269983360358
$arrBE[3,2] = 14.93410467066694
$arrBE[3,3] = 16.35856287909722
$arrBE[4,0] = 17.92637784254378
$arrBE[4,1] = 14.33155159176782
$arrBE[4,2] = 14.93227948706093
$arrBE[4,3] = 16.3568213921664
$arrBE[5,0] = 18.07684872524539
$arrBE[5,1] = 14.47807579949057
$arrBE[5,2] = 14.9452752689342
$arrBE[5,3] = 16.36926880392932
$arrBE[6,0] = 18.73645189579097
$arrBE[6,1] = 15.11649485599122
$arrBE[6,2] = 15.01067020607508
$arrBE[6,3] = 16.4330551878596
$arrBE[7,0] = 20.00081065905357
$arrBE[7,1] = 16.32578992741453
$arrBE[7,2] = 15.1673235492806
$arrBE[7,3] = 16.58951600067678
$arrBE[8,0] = 20.90118385969664
$arrBE[8,1] = 17.17786996269941
$arrBE[8,2] = 15.29919046916258
$arrBE[8,3] = 16.72308449833948
$arrBE[9,0] = 21.30256901650897
$arrBE[9,1] = 17.55573365555819
$arrBE[9,2] = 15.36268379150982
$arrBE[9,3] = 16.78775873260767
$arrBE[10,0] = 21.45323445513224
$arrBE[10,1] = 17.69728426865057
$arrBE[10,2] = 15.38721802556929
$arrBE[10,3] = 16.81279907327742
$arrBE[11,0] = 21.42084739343304
$arrBE[11,1] = 17.66686922614045
$arrBE[11,2] = 15.38191254047533
$arrBE[11,3] = 16.80738196279738
$arrBE[12,0] = 21.31499183498762
$arrBE[12,1] = 17.56741066178489
$arrBE[12,2] = 15.36469248416532
$arrBE[12,3] = 16.78980788000357
$arrBE[13,0] = 21.24997472978598
$arrBE[13,1] = 17.50628532870062
$arrBE[13,2] = 15.35420818825484
$arrBE[13,3] = 16.77911441684151
$arrBE[14,0] = 20.87477415604225
$arrBE[14,1] = 17.15296749053839
$arrBE[14,2] = 15.29511045341406
$arrBE[14,3] = 16.71893553155521
$arrBE[15,0] = 20.64238871355699
$arrBE[15,1] = 16.9336192677889
$arrBE[15,2] = 15.25974464670897
$arrBE[15,3] = 16.68301100062178
$arrBE[16,0] = 20.50796318873774
$arrBE[16,1] = 16.80654521862619
$arrBE[16,2] = 15.23973407833232
$arrBE[16,3] = 16.66271734414155
$arrBE[17,0] = 20.46232270926428
$arrBE[17,1] = 16.76336809509675
$arrBE[17,2] = 15.23301607472031
$arrBE[17,3] = 16.65591005160468
$arrBE[18,0] = 20.66720670656649
$arrBE[18,1] = 16.95706455220361
$arrBE[18,2] = 15.26347524457806
$arrBE[18,3] = 16.68679711102553
$arrBE[19,0] = 21.34612142444477
$arrBE[19,1] = 17.59666680692944
$arrBE[19,2] = 15.36973722711239
$arrBE[19,3] = 16.79495500710566
$arrBE[20,0] = 21.7820160944516
$arrBE[20,1] = 18.00566561336352
$arrBE[20,2] = 15.44203881211828
$arrBE[20,3] = 16.86883925199886
$arrBE[21,0] = 21.55013252645107
$arrBE[21,1] = 17.78824147609973
$arrBE[21,2] = 15.40319371847231
$arrBE[21,3] = 16.82911791008399
$arrBE[22,0] = 20.6559890469891
$arrBE[22,1] = 16.94646794323632
$arrBE[22,2] = 15.26178763822734
$arrBE[22,3] = 16.68508428872435
$arrBE[23,0] = 19.66305148236946
$arrBE[23,1] = 16.00434597075198
$arrBE[23,2] = 15.12195507839646
$arrBE[23,3] = 16.5438784153395
$ws.Range("B2:E25").Value = $arrBE

$arrGJ = New-Object "object[,]" 24,4
$arrGJ[0,0] = 48.42850894251271
$arrGJ[0,1] = 18.74572410737582
$arrGJ[0,2] = 25.61290236982534
$arrGJ[0,3] = 9.397903213633118
$arrGJ[1,0] = 48.09034601788373
$arrGJ[1,1] = 18.75872968850654
$arrGJ[1,2] = 25.66374274472542
$arrGJ[1,3] = 9.409937906010668
$arrGJ[2,0] = 47.89986158433803
$arrGJ[2,1] = 18.77149670012205
$arrGJ[2,2] = 25.70204080166868
$arrGJ[2,3] = 9.418890670657325
$arrGJ[3,0] = 47.82660944286967
$arrGJ[3,1] = 18.77789748566392
$arrGJ[3,2] = 25.71941837907575
$arrGJ[3,3] = 9.422931682766498
$arrGJ[4,0] = 47.814711582379
$arrGJ[4,1] = 18.779032560118
$arrGJ[4,2] = 25.72241056214501
$arrGJ[4,3] = 9.423626392761225
$arrGJ[5,0] = 47.89885590764341
$arrGJ[5,1] = 18.77157817834534
$arrGJ[5,2] = 25.70226800560899
$arrGJ[5,3] = 9.418943579830138
$arrGJ[6,0] = 48.30838968449606
$arrGJ[6,1] = 18.74921364650399
$arrGJ[6,2] = 25.62895645212679
$arrGJ[6,3] = 9.401728039080144
$arrGJ[7,0] = 49.24449182504273
$arrGJ[7,1] = 18.74346181955348
$arrGJ[7,2] = 25.54180656065944
$arrGJ[7,3] = 9.380391616336242
$arrGJ[8,0] = 50.00864052707728
$arrGJ[8,1] = 18.7626573311654
$arrGJ[8,2] = 25.51283215279351
$arrGJ[8,3] = 9.372311601270978
$arrGJ[9,0] = 50.37170404225728
$arrGJ[9,1] = 18.77650357576571
$arrGJ[9,2] = 25.50736494654825
$arrGJ[9,3] = 9.370288499462513
$arrGJ[10,0] = 50.51130180047703
$arrGJ[10,1] = 18.78248330095919
$arrGJ[10,2] = 25.50641099182225
$arrGJ[10,3] = 9.369760144895094
$arrGJ[11,0] = 50.4811447005029
$arrGJ[11,1] = 18.78116269677697
$arrGJ[11,2] = 25.50656668658867
$arrGJ[11,3] = 9.369863360529356
$arrGJ[12,0] = 50.38314711297603
$arrGJ[12,1] = 18.77698076818328
$arrGJ[12,2] = 25.50726405227742
$arrGJ[12,3] = 9.370240266174124
$arrGJ[13,0] = 50.3233927071641
$arrGJ[13,1] = 18.77451514347454
$arrGJ[13,2] = 25.50783679152544
$arrGJ[13,3] = 9.370502095900784
$arrGJ[14,0] = 49.9852152235856
$arrGJ[14,1] = 18.76185546092928
$arrGJ[14,2] = 25.51334529218439
$arrGJ[14,3] = 9.372477080956154
$arrGJ[15,0] = 49.78163587549413
$arrGJ[15,1] = 18.75539996684495
$arrGJ[15,2] = 25.51870551650849
$arrGJ[15,3] = 9.374112034868626
$arrGJ[16,0] = 49.66600332868086
$arrGJ[16,1] = 18.75216828242517
$arrGJ[16,2] = 25.52251421584708
$arrGJ[16,3] = 9.37520796009637
$arrGJ[17,0] = 49.62710625098183
$arrGJ[17,1] = 18.75115671775282
$arrGJ[17,2] = 25.52392817392529
$arrGJ[17,3] = 9.375605731241727
$arrGJ[18,0] = 49.80315681175036
$arrGJ[18,1] = 18.75603733805367
$arrGJ[18,2] = 25.51805975605379
$arrGJ[18,3] = 9.373921892013787
$arrGJ[19,0] = 50.41187490621214
$arrGJ[19,1] = 18.77818910989879
$arrGJ[19,2] = 25.50702886804131
$arrGJ[19,3] = 9.370123106860801
$arrGJ[20,0] = 50.82197291881106
$arrGJ[20,1] = 18.79695918439356
$arrGJ[20,2] = 25.50632931174261
$arrGJ[20,3] = 9.369026165922492
$arrGJ[21,0] = 50.60201052847235
$arrGJ[21,1] = 18.78654829920301
$arrGJ[21,2] = 25.50610486245861
$arrGJ[21,3] = 9.369484807571107
$arrGJ[22,0] = 49.79342280294519
$arrGJ[22,1] = 18.75574768856798
$arrGJ[22,2] = 25.51834944035669
$arrGJ[22,3] = 9.374007369801546
$arrGJ[23,0] = 48.97746509593344
$arrGJ[23,1] = 18.74091615420148
$arrGJ[23,2] = 25.55926352410763
$arrGJ[23,3] = 9.384831084233838
$ws.Range("G2:J25").Value = $arrGJ

$arrN = New-Object "object[,]" 24,1
$arrN[0,0] = 18.03386714620476
$arrN[1,0] = 18.11037131698798
$arrN[2,0] = 18.15931704480719
$arrN[3,0] = 18.17976073990558
$arrN[4,0] = 18.18318553342444
$arrN[5,0] = 18.15959073691734
$arrN[6,0] = 18.0598377811151
$arrN[7,0] = 17.87977197357248
$arrN[8,0] = 17.75681956708456
$arrN[9,0] = 17.70288424251564
$arrN[10,0] = 17.68274516930013
$arrN[11,0] = 17.68706983269737
$arrN[12,0] = 17.70122168723855
$arrN[13,0] = 17.70992717091239
$arrN[14,0] = 17.76038436718878
$arrN[15,0] = 17.79184810426811
$arrN[16,0] = 17.81013323872236
$arrN[17,0] = 17.81635662573054
$arrN[18,0] = 17.78847929042325
$arrN[19,0] = 17.69705722395401
$arrN[20,0] = 17.63896825748086
$arrN[21,0] = 17.66982014629444
$arrN[22,0] = 17.79000171985127
$arrN[23,0] = 17.92683408821942
$ws.Range("N2:N25").Value = $arrN
